$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in column A with the new topic labels, in the same order the author
# originally typed them (this matters for shared-string ordering): first
# the six new "maths" topics down rows 3-8, then x/y/z, then the later
# letters b,c,e,f,d and finally "A" (typed last, out of row order).
$ws.Cells.Item(3, 1).Value = "subtracting"
$ws.Cells.Item(4, 1).Value = "multiplying"
$ws.Cells.Item(5, 1).Value = "decimals"
$ws.Cells.Item(6, 1).Value = "fractions"
$ws.Cells.Item(7, 1).Value = "fdp"
$ws.Cells.Item(8, 1).Value = "percentages"

$ws.Cells.Item(9, 1).Value = "x"
$ws.Cells.Item(10, 1).Value = "y"
$ws.Cells.Item(11, 1).Value = "z"

$ws.Cells.Item(13, 1).Value = "b"
$ws.Cells.Item(14, 1).Value = "c"
$ws.Cells.Item(16, 1).Value = "e"
$ws.Cells.Item(17, 1).Value = "f"
$ws.Cells.Item(15, 1).Value = "d"
$ws.Cells.Item(12, 1).Value = "A"

# Fill in columns B & C for every new row (3-17): B holds the hyperlink
# text (the shared YouTube URL), C holds the "adding" exercise label.
for ($row = 3; $row -le 17; $row++) {
    $ws.Cells.Item($row, 2).Value = "https://www.youtube.com/embed/PGc33iLXaTE"
    $ws.Cells.Item($row, 3).Value = "adding"
}

# Add a single hyperlink spanning B9:B17 pointing at the YouTube embed
$ws.Hyperlinks.Add($ws.Range("B9:B17"), "https://www.youtube.com/embed/PGc33iLXaTE", [Type]::Missing, [Type]::Missing, "https://www.youtube.com/embed/PGc33iLXaTE")

# Re-apply the same "Hyperlink" cell style already used by B2:B8 so the new
# rows visually match (Hyperlinks.Add's own auto-styling only touched the
# first cell of the range).
$ws.Range("B9:B17").Style = $ws.Range("B8").Style

# Update selection / view state to match the saved file
$ws.Range("C8:C17").Select()
